$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.817
$ws.Range("D4").Value = -8.300000000000001
$ws.Range("B7").Value = 5.481
$ws.Range("A8").Value = -22.322
$ws.Range("A10").Value = -21.69
$ws.Range("E10").Value = 16.504
$ws.Range("D11").Value = -7.834000000000001
$ws.Range("A12").Value = -21.504
$ws.Range("E12").Value = 17.797
$ws.Range("E13").Value = 16.636
$ws.Range("B14").Value = 6.085
$ws.Range("D14").Value = -7.616
$ws.Range("E14").Value = 16.766
$ws.Range("B15").Value = 5.242000000000001
$ws.Range("A18").Value = -21.808
$ws.Range("B18").Value = 5.92
$ws.Range("D18").Value = -8.753000000000002
$ws.Range("D19").Value = -8.053999999999998
$ws.Range("B20").Value = 7.423
$ws.Range("D21").Value = -8.425000000000001
$ws.Range("A25").Value = -21.818
$ws.Range("D27").Value = -7.821000000000001
$ws.Range("B29").Value = 5.331
$ws.Range("E29").Value = 16.816
$ws.Range("B30").Value = 5.306
$ws.Range("B31").Value = 5.792999999999999
$ws.Range("D31").Value = -7.726999999999999
$ws.Range("E32").Value = 16.656
$ws.Range("B35").Value = 8.191999999999998
$ws.Range("E35").Value = 16.703
$ws.Range("A37").Value = -20.712
$ws.Range("D38").Value = -8.068999999999999
$ws.Range("B40").Value = 8.628
$ws.Range("D42").Value = -8.303999999999998
$ws.Range("E43").Value = 16.642
$ws.Range("B44").Value = 5.081
$ws.Range("D44").Value = -7.487
$ws.Range("D47").Value = -8.1
$ws.Range("E48").Value = 16.731
$ws.Range("E49").Value = 16.226
$ws.Range("B50").Value = 4.853
$ws.Range("E50").Value = 16.417
$ws.Range("E51").Value = 16.497
$ws.Range("B54").Value = 4.927
$ws.Range("A55").Value = -21.843
$ws.Range("D56").Value = -8.263999999999999
$ws.Range("E56").Value = 16.139
$ws.Range("D58").Value = -8.370999999999999
$ws.Range("E61").Value = 16.375
$ws.Range("D65").Value = -8.027000000000001
$ws.Range("A68").Value = -21.447
$ws.Range("B68").Value = 5.298
$ws.Range("E69").Value = 17.046
$ws.Range("E71").Value = 17.316
$ws.Range("D73").Value = -8.154
$ws.Range("B76").Value = 6.33
$ws.Range("A77").Value = -21.036
$ws.Range("A78").Value = -20.27
$ws.Range("A79").Value = -21.775
$ws.Range("E79").Value = 16.896
$ws.Range("A80").Value = -20.864
$ws.Range("A81").Value = -21.71
$ws.Range("E81").Value = 16.673
$ws.Range("A82").Value = -22.005
$ws.Range("A84").Value = -21.873
$ws.Range("B87").Value = 4.636
$ws.Range("B88").Value = 4.805
$ws.Range("D90").Value = -8.211000000000002
$ws.Range("B92").Value = 5.715
$ws.Range("D92").Value = -6.667
$ws.Range("E92").Value = 17.631
$ws.Range("D94").Value = -7.348999999999999
$ws.Range("D95").Value = -7.744000000000002
$ws.Range("B96").Value = 5.756
$ws.Range("B98").Value = 6.053
$ws.Range("A101").Value = -21.618
$ws.Range("B101").Value = 5.520999999999999
$ws.Range("D101").Value = -7.858
$ws.Range("A102").Value = -20.934
$ws.Range("B102").Value = 6.353
